$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for week 41 (bold + centered, matching the other header cells; stored as text)
$ws.Range("AR1").NumberFormat = "@"
$ws.Range("AR1").Value = "41"
$ws.Range("AR1").Font.Bold = $true
$ws.Range("AR1").HorizontalAlignment = -4108

# AR column (week 41) data values
$ws.Range("AR2").Value = 0
$ws.Range("AR5").Value = 0
$ws.Range("AR6").Value = 20
$ws.Range("AR7").Value = 2
$ws.Range("AR8").Value = 16
$ws.Range("AR9").Value = 0
$ws.Range("AR11").Value = 0
$ws.Range("AR13").Value = 0
$ws.Range("AR14").Value = 0
$ws.Range("AR16").Value = 0
$ws.Range("AR17").Value = 0
$ws.Range("AR22").Value = 0
$ws.Range("AR23").Value = 0
$ws.Range("AR24").Value = 0
$ws.Range("AR25").Value = 2
$ws.Range("AR26").Value = 0
$ws.Range("AR29").Value = 2
$ws.Range("AR30").Value = 3
$ws.Range("AR31").Value = 0
$ws.Range("AR35").Value = 0
$ws.Range("AR36").Value = 0
$ws.Range("AR37").Value = 0
$ws.Range("AR38").Value = 0
$ws.Range("AR41").Value = 0
$ws.Range("AR42").Value = 0
$ws.Range("AR43").Value = 0
$ws.Range("AQ44").Value = 0
$ws.Range("AR44").Value = 0
$ws.Range("AR45").Value = 0
$ws.Range("AR46").Value = 0
$ws.Range("AR47").Value = 0
$ws.Range("AR48").Value = 0
$ws.Range("AR49").Value = 0
$ws.Range("AR50").Value = 0
$ws.Range("AR51").Value = 0
$ws.Range("AR53").Value = 0
$ws.Range("AR54").Value = 1
$ws.Range("AR55").Value = 0
$ws.Range("AR56").Value = 0
$ws.Range("AR57").Value = 0
$ws.Range("AR58").Value = 0

# Row 36 corrections (weeks 8-14 columns K..U)
$ws.Range("K36").Value = 1
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 1
$ws.Range("Q36").Value = 2
$ws.Range("S36").Value = 3
$ws.Range("U36").Value = 0

Write-Output "done"
